# Trading update: 2026-02-17 12:48:10
#
# Appends one new trade row (row 46) to both the "All Trades" sheet and
# the "MarketMaking" sheet, mirroring the existing layout used by rows
# 2-45 of those sheets:
#   Trade # | Date | Time | Strategy | Side | Entry Price | Exit Price |
#   Status | P&L % | P&L $ | Capital After | Entry Slippage (bps) |
#   Exit Slippage (bps) | Confidence | Entry Reason | Exit Reason |
#   Duration (min)
#
# The new trade is an OPEN MarketMaking position, so Exit Price and Exit
# Reason are left blank (no exit yet).

$wb = $excel.ActiveWorkbook

$newRow = 46

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # A: Trade #
    $ws.Cells.Item($newRow, 1).Value = 45

    # B: Date - force text so Excel does not auto-convert the
    # "yyyy-mm-dd" looking string into a date serial number, then strip
    # the temporary text format back off so the cell keeps the plain,
    # unstyled look used by the rest of the sheet.
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2026-02-17"
    $ws.Cells.Item($newRow, 2).ClearFormats()

    # C: Time - same trick, to avoid "hh:mm:ss" auto-conversion.
    $ws.Cells.Item($newRow, 3).NumberFormat = "@"
    $ws.Cells.Item($newRow, 3).Value = "12:48:05"
    $ws.Cells.Item($newRow, 3).ClearFormats()

    # D: Strategy
    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"

    # E: Side
    $ws.Cells.Item($newRow, 5).Value = "UP"

    # F: Entry Price
    $ws.Cells.Item($newRow, 6).Value = 0.935975

    # G: Exit Price - left blank, trade is still OPEN.

    # H: Status
    $ws.Cells.Item($newRow, 8).Value = "OPEN"

    # I: P&L %
    $ws.Cells.Item($newRow, 9).Value = 0

    # J: P&L $
    $ws.Cells.Item($newRow, 10).Value = 0

    # K: Capital After
    $ws.Cells.Item($newRow, 11).Value = 100.118799984049

    # L: Entry Slippage (bps)
    $ws.Cells.Item($newRow, 12).Value = 0

    # M: Exit Slippage (bps)
    $ws.Cells.Item($newRow, 13).Value = 0

    # N: Confidence
    $ws.Cells.Item($newRow, 14).Value = 0.6

    # O: Entry Reason
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"

    # P: Exit Reason - left blank, trade is still OPEN.

    # Q: Duration (min)
    $ws.Cells.Item($newRow, 17).Value = 0
}
